$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set values and formulas first (so whole-column dependent formulas recalc correctly) ---

# Row 31 values (A-F)
$ws.Range("A31").Value = 45639
$ws.Range("B31").Value = "Wizards"
$ws.Range("C31").Value = "Cavaliers"
$ws.Range("D31").Value = 232.5
$ws.Range("E31").Value = "L"
$ws.Range("F31").Formula = "=IF(E31=""L"",-1,IF(E31=""W"",1/1.1,0))"

# Row 32 values (A-F)
$ws.Range("A32").Value = 45640
$ws.Range("B32").Value = "Nets"
$ws.Range("C32").Value = "Grizzlies"
$ws.Range("D32").Value = 228.5
$ws.Range("E32").Value = "W"
$ws.Range("F32").Formula = "=IF(E32=""L"",-1,IF(E32=""W"",1/1.1,0))"

# Row 35 values (J-P)
$ws.Range("J35").Value = 45639
$ws.Range("K35").Value = "Spurs"
$ws.Range("L35").Value = "Trail Blazers"
$ws.Range("M35").Value = "Spurs ML"
$ws.Range("N35").Value = -160
$ws.Range("O35").Value = "W"
$ws.Range("P35").Formula = "=IF(O35=""L"",-1,IF(N35<0,1/(-N35/100),1*(N35/100)))"

# --- Now copy formats from neighboring rows to match style indices exactly ---
$ws.Range("A30:F30").Copy()
$ws.Range("A31:F31").PasteSpecial(-4122)

$ws.Range("A30:F30").Copy()
$ws.Range("A32:F32").PasteSpecial(-4122)

$ws.Range("J34:P34").Copy()
$ws.Range("J35:P35").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("N20").Select() | Out-Null
